$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "30.066.13"
Set-TextValue $ws.Range("E2") "  +7.58%  "

Set-TextValue $ws.Range("D3") "1.884.48"
Set-TextValue $ws.Range("E3") "  +5.80%  "

Set-TextValue $ws.Range("D4") "0.9997"
Set-TextValue $ws.Range("E4") "  -0.06%  "

Set-TextValue $ws.Range("D5") "249.31"
Set-TextValue $ws.Range("E5") "  +2.16%  "

Set-TextValue $ws.Range("D6") "0.9997"
Set-TextValue $ws.Range("E6") "  -0.02%  "

Set-TextValue $ws.Range("D7") "0.4987"
Set-TextValue $ws.Range("E7") "  +1.51%  "

Set-TextValue $ws.Range("D8") "45.80"
Set-TextValue $ws.Range("E8") "  +8.94%  "

Set-TextValue $ws.Range("D9") "0.2866"
Set-TextValue $ws.Range("E9") "  +7.30%  "

Set-TextValue $ws.Range("D10") "0.06582"
Set-TextValue $ws.Range("E10") "  +5.24%  "

Set-TextValue $ws.Range("D11") "1.880.45"
Set-TextValue $ws.Range("E11") "  +5.55%  "

Set-TextValue $ws.Range("D12") "17.20"
Set-TextValue $ws.Range("E12") "  +5.28%  "

Set-TextValue $ws.Range("D13") "0.07210"
Set-TextValue $ws.Range("E13") "  +2.33%  "

Set-TextValue $ws.Range("D14") "0.6652"
Set-TextValue $ws.Range("E14") "  +6.21%  "

Set-TextValue $ws.Range("D15") "85.27"
Set-TextValue $ws.Range("E15") "  +6.57%  "

Set-TextValue $ws.Range("D16") "4.797"
Set-TextValue $ws.Range("E16") "  +3.62%  "

Set-TextValue $ws.Range("D17") "30.035.49"
Set-TextValue $ws.Range("E17") "  +7.54%  "

Set-TextValue $ws.Range("E18") "  -0.09%  "

Set-TextValue $ws.Range("D19") "12.93"
Set-TextValue $ws.Range("E19") "  +8.36%  "

Set-TextValue $ws.Range("D20") "0.000007519"
Set-TextValue $ws.Range("E20") "  +4.25%  "

Set-TextValue $ws.Range("D21") "0.9989"
Set-TextValue $ws.Range("E21") "  -0.09%  "

Set-TextValue $ws.Range("D22") "2.119.95"
Set-TextValue $ws.Range("E22") "  +5.44%  "

Set-TextValue $ws.Range("D23") "4.770"
Set-TextValue $ws.Range("E23") "  +4.17%  "

Set-TextValue $ws.Range("D24") "5.538"
Set-TextValue $ws.Range("E24") "  +5.87%  "

Set-TextValue $ws.Range("D25") "9.041"
Set-TextValue $ws.Range("E25") "  +4.03%  "

Set-TextValue $ws.Range("D26") "144.83"
Set-TextValue $ws.Range("E26") "  +2.32%  "

Set-TextValue $ws.Range("D27") "135.59"
Set-TextValue $ws.Range("E27") "  +24.08%  "

Set-TextValue $ws.Range("D28") "16.73"
Set-TextValue $ws.Range("E28") "  +6.55%  "

Set-TextValue $ws.Range("D29") "1.977"
Set-TextValue $ws.Range("E29") "  +6.23%  "

Set-TextValue $ws.Range("D30") "1.391"
Set-TextValue $ws.Range("E30") "  -0.45%  "

Set-TextValue $ws.Range("D31") "4.184"
Set-TextValue $ws.Range("E31") "  -0.50%  "

Set-TextValue $ws.Range("D32") "0.08617"
Set-TextValue $ws.Range("E32") "  +4.11%  "

Set-TextValue $ws.Range("D33") "3.891"
Set-TextValue $ws.Range("E33") "  +2.56%  "

Set-TextValue $ws.Range("D34") "0.05083"
Set-TextValue $ws.Range("E34") "  +4.56%  "

Set-TextValue $ws.Range("D35") "1.141"
Set-TextValue $ws.Range("E35") "  +6.65%  "

Set-TextValue $ws.Range("D36") "0.6877"
Set-TextValue $ws.Range("E36") "  +5.88%  "

Set-TextValue $ws.Range("D37") "1.000"
Set-TextValue $ws.Range("E37") "  +0.01%  "

Set-TextValue $ws.Range("D38") "2.709"
Set-TextValue $ws.Range("E38") "  +3.74%  "

Set-TextValue $ws.Range("D39") "2.306"
Set-TextValue $ws.Range("E39") "  +12.87%  "

Set-TextValue $ws.Range("D40") "2.749"
Set-TextValue $ws.Range("E40") "  +6.42%  "

Set-TextValue $ws.Range("D41") "0.9623"
Set-TextValue $ws.Range("E41") "  +1.66%  "

Set-TextValue $ws.Range("D42") "0.01635"
Set-TextValue $ws.Range("E42") "  +5.42%  "

Set-TextValue $ws.Range("D43") "6.095"
Set-TextValue $ws.Range("E43") "  +2.48%  "

Set-TextValue $ws.Range("E44") "  +0.07%  "

Set-TextValue $ws.Range("D45") "103.60"
Set-TextValue $ws.Range("E45") "  +3.58%  "

Set-TextValue $ws.Range("D46") "0.4213"
Set-TextValue $ws.Range("E46") "  +5.86%  "

Set-TextValue $ws.Range("D47") "7.468"
Set-TextValue $ws.Range("E47") "  +4.26%  "

Set-TextValue $ws.Range("E48") "  +4.55%  "

Set-TextValue $ws.Range("D49") "0.05635"
Set-TextValue $ws.Range("E49") "  +4.03%  "

Set-TextValue $ws.Range("D50") "32.49"
Set-TextValue $ws.Range("E50") "  +5.91%  "

Set-TextValue $ws.Range("D51") "8.306"
Set-TextValue $ws.Range("E51") "  +4.18%  "
